$d = $word.ActiveDocument

$replacements = @(
    @{old = "174÷9=19, 3"; new = "556÷7=79, 3"},
    @{old = "183÷5=36, 3"; new = "182÷7=26, 0"},
    @{old = "522÷9=58, 0"; new = "504÷9=56, 0"},
    @{old = "797÷7=113, 6"; new = "775÷9=86, 1"},
    @{old = "212÷8=26, 4"; new = "511÷3=170, 1"},
    @{old = "305÷4=76, 1"; new = "749÷5=149, 4"},
    @{old = "581÷2=290, 1"; new = "988÷5=197, 3"},
    @{old = "683÷2=341, 1"; new = "128÷9=14, 2"},
    @{old = "153÷8=19, 1"; new = "390÷7=55, 5"},
    @{old = "509÷9=56, 5"; new = "131÷9=14, 5"},
    @{old = "335÷2=167, 1"; new = "459÷7=65, 4"},
    @{old = "697÷4=174, 1"; new = "974÷6=162, 2"},
    @{old = "118÷2=59, 0"; new = "621÷2=310, 1"},
    @{old = "658÷5=131, 3"; new = "275÷8=34, 3"},
    @{old = "154÷5=30, 4"; new = "167÷8=20, 7"},
    @{old = "478÷6=79, 4"; new = "798÷7=114, 0"},
    @{old = "123÷4=30, 3"; new = "922÷8=115, 2"},
    @{old = "701÷7=100, 1"; new = "886÷2=443, 0"},
    @{old = "415÷7=59, 2"; new = "894÷9=99, 3"},
    @{old = "664÷9=73, 7"; new = "436÷2=218, 0"},
    @{old = "848÷8=106, 0"; new = "374÷3=124, 2"},
    @{old = "458÷8=57, 2"; new = "931÷3=310, 1"},
    @{old = "212÷2=106, 0"; new = "541÷3=180, 1"},
    @{old = "631÷2=315, 1"; new = "158÷8=19, 6"},
    @{old = "456÷2=228, 0"; new = "203÷8=25, 3"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
